# Delete the "datetime" column (column A) from Sheet1.
# This shifts all remaining columns (age, gender, location, mood, activity,
# period, song-artist, Id, song, artist) one position to the left, and
# removes the now-unused "datetime" shared string, the m/d/yyyy h:mm:ss
# number format, and the style that used it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns("A")
$col.Select()
$col.Delete()
